$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot the existing data rows (2..18, cols A..E) before shifting ---
$data = @{}
for ($r = 2; $r -le 18; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $data[$r,$c] = $ws.Cells.Item($r, $c).Value2
    }
}

# --- Step 2: shift rows 2..18 down into rows 3..19 (bottom-up so nothing is clobbered) ---
for ($r = 18; $r -ge 2; $r--) {
    $dest = $r + 1
    for ($c = 1; $c -le 5; $c++) {
        $v = $data[$r,$c]
        if ($v -eq $null) {
            $ws.Cells.Item($dest, $c).Value2 = $null
        } else {
            $ws.Cells.Item($dest, $c).Value2 = $v
        }
    }
}

# --- Step 3: the new last row (19) needs column A's bordered/date-formatted style,
#     which it doesn't have yet since it was previously the header-adjacent style holder.
#     Copy formats only (xlPasteFormats) from the row above so the value we already wrote stays put. ---
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Step 4: write the new row 2 (2007 data), reusing column A's pre-existing bordered style ---
$ws.Cells.Item(2, 1).Value2 = 39400
$ws.Cells.Item(2, 2).Value2 = 2007
$ws.Cells.Item(2, 3).Value2 = 3.145939949069287
$ws.Cells.Item(2, 4).Value2 = 2008
$ws.Cells.Item(2, 5).Value2 = 2.51031180018495

# --- Step 5: refresh the other rows' recomputed/simulated values per the bugfix ---
$newValues = @{
    3  = @(39765, 2008, 1.769627576887389,   2009, -1.56363396419209)
    4  = @(40130, 2009, -4.774178217057779,  2010, 1.68591416918662)
    5  = @(40494, 2010, 1.97975191822708,    2011, 3.609042024648068)
    6  = @(40862, 2011, 3.452886745653183,   2012, 2.11424984565185)
    7  = @(41228, 2012, 1.239479831392853,   2013, 0.9276342348636168)
    8  = @(41592, 2013, 0.2379616621361214,  2014, 1.3307042289459)
    9  = @(41957, 2014, 1.51977456621637,    2015, 0.3626364251072101)
    10 = @(42321, 2015, 1.470039379455756,   2016, 1.339087911421144)
    11 = @(42689, 2016, 1.638797242243251,   2017, 1.006353890555189)
    12 = @(43053, 2017, 2.161565493242668,   2018, 3.257358596620663)
    13 = @(43418, 2018, 2.214251681313772,   2019, -0.3013396321239648)
    14 = @(43783, 2019, 0.6066442151010376,  2020, 0.2691345740889695)
    15 = @(44159, 2020, -4.207901339433196,  2021, -0.6332519459683494)
    16 = @(44525, 2021, 1.099928004397532,   2022, 0.4297312830033428)
    17 = @(44890, 2022, 2.310042359896225,   2023, -0.177017417229286)
    18 = @(45254, 2023, 0.0464415346324687,  2024, 0.08221134935635366)
    19 = @(45618, 2024, -0.3101476031197148, 2025, 0.5295174046934692)
}

foreach ($r in $newValues.Keys) {
    $row = $newValues[$r]
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $row[$c - 1]
    }
}
